$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1908127208480565
$ws.Range("C2").Value = 0.558303886925795
$ws.Range("J2").Value = 0.007067137809187279
$ws.Range("P2").Value = 0.1307420494699647
$ws.Range("S2").Value = 0.1130742049469965
$ws.Range("B3").Value = 0.006289308176100629
$ws.Range("C3").Value = 0.006289308176100629
$ws.Range("J3").Value = 0.01257861635220126
$ws.Range("P3").Value = 0.7610062893081762
$ws.Range("S3").Value = 0.2138364779874214
$ws.Range("P4").Value = 0.8125
$ws.Range("S4").Value = 0.1875
$ws.Range("B6").Value = 0.08928571428571429
$ws.Range("D6").Value = 0.02380952380952381
$ws.Range("E6").Value = 0.0119047619047619
$ws.Range("F6").Value = 0.05952380952380952
$ws.Range("J6").Value = 0.1845238095238095
$ws.Range("O6").Value = 0.005952380952380952
$ws.Range("Q6").Value = 0.1547619047619048
$ws.Range("R6").Value = 0.07738095238095238
$ws.Range("S6").Value = 0.3928571428571428
$ws.Range("B7").Value = 0.08771929824561403
$ws.Range("D7").Value = 0.04093567251461988
$ws.Range("F7").Value = 0.0935672514619883
$ws.Range("J7").Value = 0.1286549707602339
$ws.Range("O7").Value = 0.01754385964912281
$ws.Range("Q7").Value = 0.1637426900584795
$ws.Range("R7").Value = 0.04678362573099415
$ws.Range("S7").Value = 0.4210526315789473
$ws.Range("B8").Value = 0.09214092140921409
$ws.Range("D8").Value = 0.02168021680216802
$ws.Range("E8").Value = 0.002710027100271003
$ws.Range("F8").Value = 0.04065040650406504
$ws.Range("J8").Value = 0.1219512195121951
$ws.Range("O8").Value = 0.01626016260162602
$ws.Range("Q8").Value = 0.1788617886178862
$ws.Range("R8").Value = 0.08401084010840108
$ws.Range("S8").Value = 0.4417344173441735
$ws.Range("B9").Value = 0.1630434782608696
$ws.Range("D9").Value = 0.04347826086956522
$ws.Range("F9").Value = 0.04347826086956522
$ws.Range("J9").Value = 0.05434782608695652
$ws.Range("O9").Value = 0.005434782608695652
$ws.Range("Q9").Value = 0.1521739130434783
$ws.Range("R9").Value = 0.1413043478260869
$ws.Range("S9").Value = 0.3967391304347826
$ws.Range("B10").Value = 0.1285574092247301
$ws.Range("D10").Value = 0.02355250245338567
$ws.Range("F10").Value = 0.07262021589793916
$ws.Range("J10").Value = 0.08832188420019627
$ws.Range("O10").Value = 0.01766437684003926
$ws.Range("Q10").Value = 0.2158979391560353
$ws.Range("R10").Value = 0.08341511285574092
$ws.Range("S10").Value = 0.3699705593719333
$ws.Range("G11").Value = 0.1397379912663755
$ws.Range("J11").Value = 0.07423580786026202
$ws.Range("K11").Value = 0.1441048034934498
$ws.Range("L11").Value = 0.6375545851528385
$ws.Range("S11").Value = 0.004366812227074236
$ws.Range("G12").Value = 0.8095238095238095
$ws.Range("J12").Value = 0.163265306122449
$ws.Range("K12").Value = 0.006802721088435374
$ws.Range("L12").Value = 0.006802721088435374
$ws.Range("S12").Value = 0.01360544217687075
$ws.Range("F15").Value = 0.02105263157894737
$ws.Range("H15").Value = 0.2210526315789474
$ws.Range("I15").Value = 0.05789473684210526
$ws.Range("J15").Value = 0.3210526315789474
$ws.Range("K15").Value = 0.07368421052631578
$ws.Range("M15").Value = 0.01052631578947368
$ws.Range("O15").Value = 0.06842105263157895
$ws.Range("S15").Value = 0.2263157894736842
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.1798941798941799
$ws.Range("I16").Value = 0.07936507936507936
$ws.Range("J16").Value = 0.417989417989418
$ws.Range("K16").Value = 0.09523809523809523
$ws.Range("M16").Value = 0.005291005291005291
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.09523809523809523
$ws.Range("S16").Value = 0.1058201058201058
$ws.Range("F17").Value = 0.01095890410958904
$ws.Range("H17").Value = 0.189041095890411
$ws.Range("I17").Value = 0.1287671232876712
$ws.Range("J17").Value = 0.3917808219178082
$ws.Range("K17").Value = 0.0958904109589041
$ws.Range("M17").Value = 0.02191780821917808
$ws.Range("O17").Value = 0.06027397260273973
$ws.Range("S17").Value = 0.1013698630136986
$ws.Range("F18").Value = 0.006097560975609756
$ws.Range("H18").Value = 0.1524390243902439
$ws.Range("I18").Value = 0.1036585365853658
$ws.Range("J18").Value = 0.451219512195122
$ws.Range("K18").Value = 0.09146341463414634
$ws.Range("M18").Value = 0.01829268292682927
$ws.Range("O18").Value = 0.04878048780487805
$ws.Range("S18").Value = 0.1280487804878049
$ws.Range("F19").Value = 0.01083743842364532
$ws.Range("H19").Value = 0.2019704433497537
$ws.Range("I19").Value = 0.09655172413793103
$ws.Range("J19").Value = 0.4098522167487685
$ws.Range("K19").Value = 0.1064039408866995
$ws.Range("M19").Value = 0.02167487684729064
$ws.Range("N19").Value = 0.0009852216748768472
$ws.Range("O19").Value = 0.0748768472906404
$ws.Range("S19").Value = 0.07684729064039408
